$wb = $excel.ActiveWorkbook

# Mapping of worksheet name -> hashtable of row -> new F value
$changes = @{
    "展览" = @{
        3  = 51
        4  = 73
        6  = 32
        7  = 2625
        8  = 1140
        9  = 224
        10 = 88
        11 = 5785
        13 = 228
        14 = 574
        15 = 11545
        16 = 11700
        18 = 71
        22 = 28
    }
    "全部类型" = @{
        3  = 51
        4  = 73
        6  = 32
        7  = 2625
        9  = 1140
        10 = 224
        11 = 88
        12 = 5785
        14 = 228
        15 = 574
        16 = 11545
        17 = 11700
        19 = 71
        23 = 28
    }
}

foreach ($sheetName in $changes.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rowMap = $changes[$sheetName]
    foreach ($row in $rowMap.Keys) {
        $ws.Cells.Item($row, 6).Value = $rowMap[$row]
    }
}
